# Append a new row (row 98) of data to each of the four worksheets,
# mirroring the existing row-97 layout/format (date-formatted A column,
# text columns B-E, numeric columns F-I).

$wb = $excel.ActiveWorkbook

function Add-Row98 {
    param(
        [string]$SheetName,
        [double]$AVal,
        [string]$BVal,
        [string]$CVal,
        [string]$DVal,
        [string]$EVal,
        [double]$FVal,
        [string]$GVal,
        [double]$HVal,
        [double]$IVal
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $r = 98

    $ws.Cells.Item($r, 1).Value = $AVal
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $BVal
    $ws.Cells.Item($r, 3).Value = $CVal
    $ws.Cells.Item($r, 4).Value = $DVal
    $ws.Cells.Item($r, 5).Value = $EVal

    $ws.Cells.Item($r, 6).Value = $FVal
    $ws.Cells.Item($r, 7).Value = [double]$GVal
    $ws.Cells.Item($r, 8).Value = $HVal
    $ws.Cells.Item($r, 9).Value = $IVal
}

# FE_LFT_#1
Add-Row98 `
    "FE_LFT_#1" `
    45884.49731481481 `
    "0x01,0x7c" `
    "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," `
    "0x01,0x00" `
    "0xf" `
    380 `
    "7.598631275147109e+23" `
    256 `
    15

# FE_LFT_#2
Add-Row98 `
    "FE_LFT_#2" `
    45884.49731481481 `
    "0x01,0x90" `
    "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," `
    "0x01,0x14" `
    "0xe" `
    400 `
    "5.68432987514711e+23" `
    276 `
    14

# FE_PLT_#1
Add-Row98 `
    "FE_PLT_#1" `
    45884.49731481481 `
    "0x00,0x6e" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," `
    "0x00,0x5C" `
    "0x3" `
    110 `
    "5.68631262647114e+23" `
    92 `
    3

# FE_PLT_#2
Add-Row98 `
    "FE_PLT_#2" `
    45884.49731481481 `
    "0x00,0x6e" `
    "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," `
    "0x00,0x5A" `
    "0x3" `
    110 `
    "9.85046333984776e+23" `
    90 `
    3

Write-Output "Row 98 added to all four sheets."
